$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6, 1).Value = "25-03-2025"
$ws.Cells.Item(6, 2).Value = "Gujarat Titans vs Punjab Kings"
$ws.Cells.Item(6, 3).Value = "Gujarat Titans"
$ws.Cells.Item(6, 4).Value = "Gujarat Titans"
